$d = $word.ActiveDocument

# Commit message: "minor wording change in appearance of graphs"
# change: 'graphs look more professional and ready-to-be-published look'
# to 'graphs have a more professional and ready-to-be-published look'
#
# Full sentence before: "...the difference is not very significant but
# ggplot2 graphs look more professional and ready-to-be-published look. "
# Full sentence after:  "...the difference is not very significant but
# ggplot2 graphs have a more professional and ready-to-be-published look. "
#
# The original text is split around a (zero-length) "_GoBack" bookmark:
#   ["the difference is not very significant but "] <bookmark/>
#   ["ggplot2 graphs look more professional and ready-to-be-published look. "]
# We insert new text just ahead of the bookmark and trim the run that
# follows it, preserving the bookmark's position in the sentence.

$bm = $d.Bookmarks("_GoBack")
$insertionPoint = $d.Range($bm.Range.Start, $bm.Range.Start)
$insertionPoint.InsertBefore("ggplot2 graphs have a ")

$bm = $d.Bookmarks("_GoBack")
$afterBookmark = $d.Range($bm.Range.End, $d.Content.End)
$afterBookmark.Find.Execute(
    "ggplot2 graphs look more professional",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "more professional",
    2)
